$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = 23

$ws.Range("B24").Value = 0.62121527777777774
$ws.Range("B24").NumberFormat = "h:mm:ss"

$ws.Range("C24").Formula = "=SUM(B2:B24)+1.2708333333"
$ws.Range("C24").NumberFormat = "[h]:mm:ss"

$ws.Range("C25").Select()
